$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50: copy the row-51 formatting down onto row 50 for the columns ---
# --- whose style actually changes (D/E/N/P already match; Q/R/S unstyled) ---
$fmtCols = @("A","B","C","F","G","H","I","J","K","L","M","O")
foreach ($col in $fmtCols) {
    $ws.Range($col + "51").Copy()
    $ws.Range($col + "50").PasteSpecial(-4122)
}

# --- Row 50: new BOM line values (13x Adhesive Thermal Pad for TO-220) ---
$ws.Range("A50").Value = 13
$ws.Range("B50").Value = "Thermal pads"
$ws.Range("C50").Value = "Thermal pad"
$ws.Range("D50").Value = "Adhesive Thermal Pad for TO-220"
$ws.Range("G50").Value = "Aavid"
$ws.Range("H50").Value = "53-77-9ACG"
$ws.Range("I50").Value = "53-77-9ACG-ND"
$ws.Range("J50").Value = "532-53-77-9ACG"
$ws.Range("K50").Value = 0.55
$ws.Range("L50").Value = 0.673

# --- Row 50: per-row cost/lookup formulas (mirrors the pattern used by ---
# --- the surrounding rows, e.g. row 43/44/45/51) ---
$ws.Range("M50").Formula = "=K50*A50"
$ws.Range("N50").Formula = "=L50*A50"
$ws.Range("P50").Formula = '=IF(NOT(I50=""),A50&","&I50,"")'
$ws.Range("Q50").Formula = '=A50&"x "&C50'
$ws.Range("R50").Formula = '=IF(NOT(J50=""),J50&"|"&A50,"")'
$ws.Range("S50").Formula = '=H50&" "&A50'

# Row height was re-applied (same 16.5pt) which flips on customHeight in the
# saved XML, matching the author's re-save of that row.
$ws.Rows.Item(50).RowHeight = 16.5

# --- Selection / scroll position left by the editor after the change ---
[void]$ws.Activate()
$w = $excel.ActiveWindow
$w.ScrollRow = 31
$w.ScrollColumn = 1
$ws.Rows.Item(50).Select() | Out-Null

Write-Host "done"
